$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.4139986809832124
$ws.Cells.Item(2, 4).Value = 0.04596563793990782
$ws.Cells.Item(2, 5).Value = 0.1749602748229151
$ws.Cells.Item(2, 6).Value = 0.996870010770678
$ws.Cells.Item(2, 7).Value = 0.002464209536163677
$ws.Cells.Item(2, 11).Value = 2.000725455984195
$ws.Cells.Item(2, 12).Value = 0.1527498577433519
$ws.Cells.Item(2, 15).Value = 3.561611638974739

$ws.Cells.Item(3, 3).Value = 0.4051178511413411
$ws.Cells.Item(3, 4).Value = 0.04481902887347644
$ws.Cells.Item(3, 5).Value = 0.1708091561282927
$ws.Cells.Item(3, 6).Value = 1.003477975750442
$ws.Cells.Item(3, 7).Value = 0.0024677502660045
$ws.Cells.Item(3, 11).Value = 1.788415386103793
$ws.Cells.Item(3, 12).Value = 0.1487713566611859
$ws.Cells.Item(3, 15).Value = 3.602852058618993

$ws.Cells.Item(4, 3).Value = 0.3998931948320603
$ws.Cells.Item(4, 4).Value = 0.04411986742900353
$ws.Cells.Item(4, 5).Value = 0.1683596402831817
$ws.Cells.Item(4, 6).Value = 1.008412737056688
$ws.Cells.Item(4, 7).Value = 0.002470037784509738
$ws.Cells.Item(4, 11).Value = 1.657860822143164
$ws.Cells.Item(4, 12).Value = 0.1464175458867061
$ws.Cells.Item(4, 15).Value = 3.631272634130227

$ws.Cells.Item(5, 3).Value = 0.3978214668213411
$ws.Cells.Item(5, 4).Value = 0.04383620304147229
$ws.Cells.Item(5, 5).Value = 0.1673863829903937
$ws.Cells.Item(5, 6).Value = 1.010643693257144
$ws.Cells.Item(5, 7).Value = 0.00247099859562998
$ws.Cells.Item(5, 11).Value = 1.604612395996014
$ws.Cells.Item(5, 12).Value = 0.1454806914137947
$ws.Cells.Item(5, 15).Value = 3.643631216280539

$ws.Cells.Item(6, 3).Value = 0.3974809218557027
$ws.Cells.Item(6, 4).Value = 0.0437891770713037
$ws.Cells.Item(6, 5).Value = 0.1672262799426498
$ws.Cells.Item(6, 6).Value = 1.011027410617288
$ws.Cells.Item(6, 7).Value = 0.002471159869347427
$ws.Cells.Item(6, 11).Value = 1.595767815648685
$ws.Cells.Item(6, 12).Value = 0.1453264757723289
$ws.Cells.Item(6, 15).Value = 3.645730210443418

$ws.Cells.Item(7, 3).Value = 0.3998650225560709
$ws.Cells.Item(7, 4).Value = 0.04411603673612774
$ws.Cells.Item(7, 5).Value = 0.1683464136462227
$ws.Cells.Item(7, 6).Value = 1.008441934590145
$ws.Cells.Item(7, 7).Value = 0.002470050626134021
$ws.Cells.Item(7, 11).Value = 1.657142879378284
$ws.Cells.Item(7, 12).Value = 0.1464048207255573
$ws.Cells.Item(7, 15).Value = 3.631436163385985

$ws.Cells.Item(8, 3).Value = 0.410889187738519
$ws.Cells.Item(8, 4).Value = 0.04556929390411568
$ws.Cells.Item(8, 5).Value = 0.1735083440453593
$ws.Cells.Item(8, 6).Value = 0.9989659373338782
$ws.Cells.Item(8, 7).Value = 0.002465406877625336
$ws.Cells.Item(8, 11).Value = 1.927563042905831
$ws.Cells.Item(8, 12).Value = 0.1513595726019901
$ws.Cells.Item(8, 15).Value = 3.575186827951626

$ws.Cells.Item(9, 3).Value = 0.4343210761982732
$ws.Cells.Item(9, 4).Value = 0.04845669077908354
$ws.Cells.Item(9, 5).Value = 0.1844206565040025
$ws.Cells.Item(9, 6).Value = 0.9873753370555605
$ws.Cells.Item(9, 7).Value = 0.00245719694332058
$ws.Cells.Item(9, 11).Value = 2.456214352125755
$ws.Cells.Item(9, 12).Value = 0.1617845040393746
$ws.Cells.Item(9, 15).Value = 3.48957181864418

$ws.Cells.Item(10, 3).Value = 0.4526484649215945
$ws.Cells.Item(10, 4).Value = 0.05059988414623007
$ws.Cells.Item(10, 5).Value = 0.1929232792678661
$ws.Cells.Item(10, 6).Value = 0.9831643943366828
$ws.Cells.Item(10, 7).Value = 0.002451705841230303
$ws.Cells.Item(10, 11).Value = 2.843531361937721
$ws.Cells.Item(10, 12).Value = 0.1698801903928597
$ws.Cells.Item(10, 15).Value = 3.44186472000365

$ws.Cells.Item(11, 3).Value = 0.4612290879202021
$ws.Cells.Item(11, 4).Value = 0.05157941815416933
$ws.Cells.Item(11, 5).Value = 0.1968976477203555
$ws.Cells.Item(11, 6).Value = 0.9821926324683119
$ws.Cells.Item(11, 7).Value = 0.002449323983872219
$ws.Cells.Item(11, 11).Value = 3.019482673451648
$ws.Cells.Item(11, 12).Value = 0.1736589436436446
$ws.Cells.Item(11, 15).Value = 3.423491686338224

$ws.Cells.Item(12, 3).Value = 0.4645134255538039
$ws.Cells.Item(12, 4).Value = 0.05195098093135897
$ws.Cells.Item(12, 5).Value = 0.1984180020261661
$ws.Cells.Item(12, 6).Value = 0.9819610987787399
$ws.Cells.Item(12, 7).Value = 0.002448438634509003
$ws.Cells.Item(12, 11).Value = 3.086074203756425
$ws.Cells.Item(12, 12).Value = 0.1751037252383298
$ws.Cells.Item(12, 15).Value = 3.417015545861005

$ws.Cells.Item(13, 3).Value = 0.4638045255078964
$ws.Cells.Item(13, 4).Value = 0.05187093040495938
$ws.Cells.Item(13, 5).Value = 0.1980898833829912
$ws.Cells.Item(13, 6).Value = 0.9820048833282868
$ws.Cells.Item(13, 7).Value = 0.002448628572874624
$ws.Cells.Item(13, 11).Value = 3.071734233164307
$ws.Cells.Item(13, 12).Value = 0.1747919489395287
$ws.Cells.Item(13, 15).Value = 3.418388849934274

$ws.Cells.Item(14, 3).Value = 0.4614985898201382
$ws.Cells.Item(14, 4).Value = 0.05160997426966674
$ws.Cells.Item(14, 5).Value = 0.1970224203982909
$ws.Cells.Item(14, 6).Value = 0.9821708451693496
$ws.Cells.Item(14, 7).Value = 0.002449250813261505
$ws.Cells.Item(14, 11).Value = 3.024961964015915
$ws.Cells.Item(14, 12).Value = 0.173777528851673
$ws.Cells.Item(14, 15).Value = 3.422949228921794

$ws.Cells.Item(15, 3).Value = 0.4600907023713319
$ws.Cells.Item(15, 4).Value = 0.05145021301225228
$ws.Cells.Item(15, 5).Value = 0.1963705682180787
$ws.Cells.Item(15, 6).Value = 0.9822902933191813
$ws.Cells.Item(15, 7).Value = 0.002449634113732138
$ws.Cells.Item(15, 11).Value = 2.996307646109017
$ws.Cells.Item(15, 12).Value = 0.1731579723200127
$ws.Cells.Item(15, 15).Value = 3.425805353120751

$ws.Cells.Item(16, 3).Value = 0.452092602362967
$ws.Cells.Item(16, 4).Value = 0.05053595941425471
$ws.Cells.Item(16, 5).Value = 0.1926656906314577
$ws.Cells.Item(16, 6).Value = 0.983246955807445
$ws.Cells.Item(16, 7).Value = 0.002451863829898588
$ws.Cells.Item(16, 11).Value = 2.832027442344383
$ws.Cells.Item(16, 12).Value = 0.169635175461849
$ws.Cells.Item(16, 15).Value = 3.443132654328224

$ws.Cells.Item(17, 3).Value = 0.4472483967369101
$ws.Cells.Item(17, 4).Value = 0.04997625080464019
$ws.Cells.Item(17, 5).Value = 0.190420165739333
$ws.Cells.Item(17, 6).Value = 0.9840760919148011
$ws.Cells.Item(17, 7).Value = 0.002453261358337537
$ws.Cells.Item(17, 11).Value = 2.731183131683508
$ws.Cells.Item(17, 12).Value = 0.1674986710766575
$ws.Cells.Item(17, 15).Value = 3.454616919234866

$ws.Cells.Item(18, 3).Value = 0.4444850453639333
$ws.Cells.Item(18, 4).Value = 0.04965475416742748
$ws.Cells.Item(18, 5).Value = 0.1891386197186051
$ws.Cells.Item(18, 6).Value = 0.9846417818197324
$ws.Cells.Item(18, 7).Value = 0.002454076109799231
$ws.Cells.Item(18, 11).Value = 2.673157614396303
$ws.Cells.Item(18, 12).Value = 0.1662788419460099
$ws.Cells.Item(18, 15).Value = 3.461535565744157

$ws.Cells.Item(19, 3).Value = 0.4435533552481274
$ws.Cells.Item(19, 4).Value = 0.0495459760572885
$ws.Cells.Item(19, 5).Value = 0.1887064302299422
$ws.Cells.Item(19, 6).Value = 0.9848485426626752
$ws.Cells.Item(19, 7).Value = 0.002454353850505853
$ws.Cells.Item(19, 11).Value = 2.653507384812428
$ws.Cells.Item(19, 12).Value = 0.1658673785356797
$ws.Cells.Item(19, 15).Value = 3.463931814071401

$ws.Cells.Item(20, 3).Value = 0.4477616995567075
$ws.Cells.Item(20, 4).Value = 0.05003578811447795
$ws.Cells.Item(20, 5).Value = 0.1906581683620985
$ws.Cells.Item(20, 6).Value = 0.9839786339409926
$ws.Cells.Item(20, 7).Value = 0.002453111458342161
$ws.Cells.Item(20, 11).Value = 2.74192053321633
$ws.Cells.Item(20, 12).Value = 0.1677251705321225
$ws.Cells.Item(20, 15).Value = 3.453361965658303

$ws.Cells.Item(21, 3).Value = 0.4621749476648347
$ws.Cells.Item(21, 4).Value = 0.05168660638507561
$ws.Cells.Item(21, 5).Value = 0.1973355433488706
$ws.Cells.Item(21, 6).Value = 0.9821183890662297
$ws.Cells.Item(21, 7).Value = 0.002449067595981726
$ws.Cells.Item(21, 11).Value = 3.038701150151041
$ws.Cells.Item(21, 12).Value = 0.1740751123626723
$ws.Cells.Item(21, 15).Value = 3.421596651213946

$ws.Cells.Item(22, 3).Value = 0.4717991100726806
$ws.Cells.Item(22, 4).Value = 0.05276919903282362
$ws.Cells.Item(22, 5).Value = 0.2017890662341344
$ws.Cells.Item(22, 6).Value = 0.9816982949388802
$ws.Cells.Item(22, 7).Value = 0.002446521466519157
$ws.Cells.Item(22, 11).Value = 3.232445138569801
$ws.Cells.Item(22, 12).Value = 0.1783059158836551
$ws.Cells.Item(22, 15).Value = 3.403642942633581

$ws.Cells.Item(23, 3).Value = 0.4666438085549487
$ws.Cells.Item(23, 4).Value = 0.05219106952932151
$ws.Cells.Item(23, 5).Value = 0.1994039386891941
$ws.Cells.Item(23, 6).Value = 0.9818494600813636
$ws.Cells.Item(23, 7).Value = 0.00244787155662951
$ws.Cells.Item(23, 11).Value = 3.12906128180731
$ws.Cells.Item(23, 12).Value = 0.1760404520117476
$ws.Cells.Item(23, 15).Value = 3.412967496261899

$ws.Cells.Item(24, 3).Value = 0.4475295679307862
$ws.Cells.Item(24, 4).Value = 0.05000887039919633
$ws.Cells.Item(24, 5).Value = 0.1905505379882229
$ws.Cells.Item(24, 6).Value = 0.9840224174822794
$ws.Cells.Item(24, 7).Value = 0.002453179192966841
$ws.Cells.Item(24, 11).Value = 2.737066305733038
$ws.Cells.Item(24, 12).Value = 0.1676227437329345
$ws.Cells.Item(24, 15).Value = 3.45392834541866

$ws.Cells.Item(25, 3).Value = 0.4277873234759966
$ws.Cells.Item(25, 4).Value = 0.04767166635276965
$ws.Cells.Item(25, 5).Value = 0.1813836379002751
$ws.Cells.Item(25, 6).Value = 0.9897580118543772
$ws.Cells.Item(25, 7).Value = 0.002459322574913475
$ws.Cells.Item(25, 11).Value = 2.31338464629323
$ws.Cells.Item(25, 12).Value = 0.158887956997404
$ws.Cells.Item(25, 15).Value = 3.510074999048754

Write-Host "Applied 380 kV case update"
